$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.754.13"
$ws.Range("E2").Value = "  +1.51%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.769.92"
$ws.Range("E3").Value = "  -1.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.06"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.63"
$ws.Range("E6").Value = "  +0.05%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.766.89"
$ws.Range("E7").Value = "  -1.07%  "

$ws.Range("E8").Value = "  -0.06%  "

$ws.Range("E10").Value = "  +3.57%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.35"
$ws.Range("E11").Value = "  +2.80%  "

$ws.Range("E12").Value = "  -0.80%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.46"
$ws.Range("E13").Value = "  -0.77%  "

$ws.Range("E14").Value = "  +0.80%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.400.34"
$ws.Range("E15").Value = "  -0.91%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.760.45"
$ws.Range("E16").Value = "  -1.25%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.812.00"
$ws.Range("E17").Value = "  +1.51%  "

$ws.Range("E18").Value = "  +0.95%  "

$ws.Range("E19").Value = "  -0.06%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.22"
$ws.Range("E20").Value = "  -1.49%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.81"
$ws.Range("E21").Value = "  +13.71%  "

$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "496.33"
$ws.Range("E22").Value = "  +0.20%  "

$ws.Range("E23").Value = "  -1.39%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.77"
$ws.Range("E24").Value = "  -0.19%  "

$ws.Range("E25").Value = "  -0.51%  "

$ws.Range("E26").Value = "  -2.78%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.42"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.25"
$ws.Range("E28").Value = "  +1.14%  "

$ws.Range("E29").Value = "  -0.17%  "

$ws.Range("E30").Value = "  +3.96%  "

$ws.Range("E31").Value = "  -0.06%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.00"
$ws.Range("E32").Value = "  +1.22%  "

$ws.Range("E33").Value = "  -2.27%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.917.04"
$ws.Range("E34").Value = "  -0.87%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.704.11"
$ws.Range("E35").Value = "  -1.08%  "

$ws.Range("E36").Value = "  -1.00%  "

$ws.Range("E37").Value = "  +0.11%  "

$ws.Range("E38").Value = "  -0.55%  "

$ws.Range("E39").Value = "  +0.36%  "

$ws.Range("E40").Value = "  -0.10%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.327"
$ws.Range("E41").Value = "  -1.04%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "441.17"
$ws.Range("E42").Value = "  -4.97%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "48.92"
$ws.Range("E43").Value = "  -0.56%  "

$ws.Range("E45").Value = "  +0.74%  "

$ws.Range("E46").Value = "  +0.87%  "

$ws.Range("E47").Value = "  +0.02%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.85"
$ws.Range("E48").Value = "  +0.72%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.835.32"
$ws.Range("E49").Value = "  -0.49%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "141.14"
$ws.Range("E50").Value = "  +0.73%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0356"
$ws.Range("E51").Value = "  +0.80%  "
